$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.971.23"
$ws.Range("E2").Value = "  -0.97%  "

$ws.Range("D3").Value = "3.414.59"
$ws.Range("E3").Value = "  -0.84%  "

$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "572.05"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.27%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "141.92"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.55%  "

$ws.Range("E7").Value = "  +0.05%  "

$ws.Range("D8").Value = "3.415.81"
$ws.Range("E8").Value = "  -0.83%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.476"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.81%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.57"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.93%  "

$ws.Range("E11").Value = "  +1.14%  "

$ws.Range("E12").Value = "  +0.77%  "

$ws.Range("D13").Value = "3.997.32"
$ws.Range("E13").Value = "  -0.85%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "28.21"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.67%  "

$ws.Range("E15").Value = "  +0.55%  "

$ws.Range("E16").Value = "  -1.14%  "

$ws.Range("D17").Value = "3.419.66"
$ws.Range("E17").Value = "  -0.76%  "

$ws.Range("D18").Value = "61.106.69"
$ws.Range("E18").Value = "  -0.90%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.31"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.35%  "

$ws.Range("E20").Value = "  +2.29%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "9.36"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.26%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "390.40"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.78%  "

$ws.Range("E23").Value = "  +0.82%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "72.78"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.70%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.997"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.48%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0000124"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.05%  "

$ws.Range("D27").Value = "3.551.74"
$ws.Range("E27").Value = "  -1.04%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.180"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.52%  "

$ws.Range("B29").Value = "Binance-PegBSC-USD"
$ws.Range("C29").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.999"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.24%  "

$ws.Range("B30").Value = "RenderToken"
$ws.Range("C30").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.43"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -5.27%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.19"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.20%  "

$ws.Range("B32").Value = "Fetch.AI"
$ws.Range("C32").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.45"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -8.14%  "

$ws.Range("B33").Value = "PancakeSwap"
$ws.Range("C33").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.17"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.12%  "

$ws.Range("E34").Value = "  -0.03%  "

$ws.Range("E35").Value = "  -0.97%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "7.01"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.03%  "

$ws.Range("D37").Value = "3.442.11"
$ws.Range("E37").Value = "  -0.61%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.13"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.92%  "

$ws.Range("B39").Value = "Monero"
$ws.Range("C39").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "167.85"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.04%  "

$ws.Range("B40").Value = "ImmutableX"
$ws.Range("C40").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.55"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.81%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0784"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.70%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "27.08"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +4.65%  "

$ws.Range("E43").Value = "  +0.61%  "

$ws.Range("E44").Value = "  +0.50%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.00"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.04%  "

$ws.Range("E46").Value = "  -0.87%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "41.86"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.40%  "

$ws.Range("D48").Value = "2.608.73"
$ws.Range("E48").Value = "  -0.78%  "

$ws.Range("E49").Value = "  -2.91%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "6.96"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.51%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "22.94"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.64%  "
